# Auto-generated update of cryptos table (rows 2-51, columns B-E)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.567.92"
$ws.Range("E2").Value = "  +0.06%  "
$ws.Range("D3").Value = "3.841.09"
$ws.Range("E3").Value = "  +0.70%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'707.56"
$ws.Range("E5").Value = "  +0.67%  "
$ws.Range("D6").Value = "'173.25"
$ws.Range("E6").Value = "  -0.73%  "
$ws.Range("D7").Value = "3.839.91"
$ws.Range("E7").Value = "  +0.75%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").Value = "'0.528"
$ws.Range("E9").Value = "  -0.45%  "
$ws.Range("E10").Value = "  +0.16%  "
$ws.Range("D11").Value = "'7.34"
$ws.Range("E11").Value = "  +0.84%  "
$ws.Range("E12").Value = "  -0.19%  "
$ws.Range("E13").Value = "  -0.43%  "
$ws.Range("D14").Value = "'37.14"
$ws.Range("E14").Value = "  +2.04%  "
$ws.Range("D15").Value = "4.486.15"
$ws.Range("E15").Value = "  +0.77%  "
$ws.Range("D16").Value = "3.865.91"
$ws.Range("E16").Value = "  +1.52%  "
$ws.Range("D17").Value = "71.475.32"
$ws.Range("E17").Value = "  +0.12%  "
$ws.Range("E18").Value = "  +0.82%  "
$ws.Range("B19").Value = "TRON"
$ws.Range("C19").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D19").Value = "'0.115"
$ws.Range("E19").Value = "  +0.44%  "
$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D20").Value = "'17.49"
$ws.Range("E20").Value = "  -1.33%  "
$ws.Range("D21").Value = "'499.70"
$ws.Range("E21").Value = "  +3.46%  "
$ws.Range("D22").Value = "'10.79"
$ws.Range("E22").Value = "  -1.38%  "
$ws.Range("D23").Value = "'0.736"
$ws.Range("E23").Value = "  +3.01%  "
$ws.Range("D24").Value = "'85.57"
$ws.Range("E24").Value = "  +1.10%  "
$ws.Range("D25").Value = "'0.0000147"
$ws.Range("E25").Value = "  +2.01%  "
$ws.Range("D26").Value = "'10.72"
$ws.Range("E26").Value = "  +1.05%  "
$ws.Range("D27").Value = "'12.24"
$ws.Range("E27").Value = "  -0.95%  "
$ws.Range("D28").Value = "3.993.91"
$ws.Range("E28").Value = "  +0.73%  "
$ws.Range("E29").Value = "  -2.38%  "
$ws.Range("B30").Value = "Dai"
$ws.Range("C30").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D30").Value = "'1.00"
$ws.Range("E30").Value = "  -0.09%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "'3.14"
$ws.Range("E31").Value = "  -0.78%  "
$ws.Range("D32").Value = "'7.53"
$ws.Range("E32").Value = "  -1.41%  "
$ws.Range("E33").Value = "  -2.58%  "
$ws.Range("E34").Value = "  -0.10%  "
$ws.Range("D35").Value = "'0.180"
$ws.Range("E35").Value = "  -4.33%  "
$ws.Range("D36").Value = "'9.27"
$ws.Range("E36").Value = "  -0.38%  "
$ws.Range("D37").Value = "3.803.88"
$ws.Range("E37").Value = "  +1.10%  "
$ws.Range("D38").Value = "'0.995"
$ws.Range("E38").Value = "  -0.67%  "
$ws.Range("E39").Value = "  +0.15%  "
$ws.Range("E40").Value = "  -2.38%  "
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").Value = "'2.32"
$ws.Range("E41").Value = "  +0.71%  "
$ws.Range("E42").Value = "  +4.84%  "
$ws.Range("B43").Value = "Filecoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D43").Value = "'6.05"
$ws.Range("E43").Value = "  +0.12%  "
$ws.Range("E45").Value = "  +0.28%  "
$ws.Range("E46").Value = "  +0.58%  "
$ws.Range("D47").Value = "'164.31"
$ws.Range("E47").Value = "  -0.56%  "
$ws.Range("D48").Value = "'431.24"
$ws.Range("E48").Value = "  +3.22%  "
$ws.Range("D49").Value = "'49.10"
$ws.Range("E49").Value = "  +0.61%  "
$ws.Range("D50").Value = "'8.76"
$ws.Range("E50").Value = "  +1.40%  "
$ws.Range("E51").Value = "  -0.13%  "
